# chore: update Sheets via scheduled runner
# Refresh currentAveragePrice*/LevePrice*/LeveProfit* market-data columns
# (H:N) for a batch of leves across the crafting-class sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2944
$ws.Range("J32").Value = 2749.5
$ws.Range("L32").Value = 2749.5
$ws.Range("N32").Value = -3401.5
$ws.Range("H53").Value = 370.625
$ws.Range("I53").Value = 324
$ws.Range("K53").Value = 324
$ws.Range("M53").Value = 313
$ws.Range("H70").Value = 5701.7036
$ws.Range("J70").Value = 6568.3076
$ws.Range("L70").Value = 19704.9228
$ws.Range("N70").Value = -20244.9228
$ws.Range("H73").Value = 5701.7036
$ws.Range("J73").Value = 6568.3076
$ws.Range("L73").Value = 19704.9228
$ws.Range("N73").Value = -21576.9228
$ws.Range("H141").Value = 2494
$ws.Range("I141").Value = 2409.9
$ws.Range("K141").Value = 7229.700000000001
$ws.Range("M141").Value = -2049.700000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H57").Value = 15000
$ws.Range("I57").Value = 15000
$ws.Range("K57").Value = 15000
$ws.Range("M57").Value = -14516
$ws.Range("H103").Value = 180000
$ws.Range("J103").Value = 180000
$ws.Range("L103").Value = 180000
$ws.Range("N103").Value = -182344

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3146.0833
$ws.Range("I105").Value = 2822.75
$ws.Range("K105").Value = 2822.75
$ws.Range("M105").Value = -1075.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 800500.2
$ws.Range("I6").Value = 800500.2
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 800500.2
$ws.Range("L6").Value = 0
$ws.Range("M6").ClearContents()
$ws.Range("N6").Value = -800387.2
$ws.Range("H7").Value = 115.36842
$ws.Range("I7").Value = 121.57143
$ws.Range("J7").Value = 98
$ws.Range("K7").Value = 121.57143
$ws.Range("L7").Value = 98
$ws.Range("M7").Value = -8.571430000000007
$ws.Range("N7").Value = -324
$ws.Range("H11").Value = 237
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 237
$ws.Range("K11").Value = 0
$ws.Range("L11").ClearContents()
$ws.Range("M11").Value = 237
$ws.Range("N11").Value = -517
$ws.Range("H12").Value = 449.5
$ws.Range("I12").Value = 449.5
$ws.Range("K12").Value = 449.5
$ws.Range("M12").Value = -279.5
$ws.Range("H16").Value = 1650.1333
$ws.Range("I16").Value = 1286.6923
$ws.Range("K16").Value = 1286.6923
$ws.Range("M16").Value = -999.6922999999999
$ws.Range("H19").Value = 231.66667
$ws.Range("I19").Value = 231.66667
$ws.Range("K19").Value = 231.66667
$ws.Range("M19").Value = -61.66667000000001
$ws.Range("H24").Value = 231.66667
$ws.Range("I24").Value = 231.66667
$ws.Range("K24").Value = 231.66667
$ws.Range("M24").Value = -61.66667000000001
$ws.Range("H58").Value = 3342.2222
$ws.Range("I58").Value = 3297.2856
$ws.Range("K58").Value = 3297.2856
$ws.Range("M58").Value = -3094.2856
$ws.Range("J62").Value = 2500
$ws.Range("L62").Value = 2500
$ws.Range("N62").Value = -3748
$ws.Range("J65").Value = 2500
$ws.Range("L65").Value = 12500
$ws.Range("N65").Value = -18740
$ws.Range("H99").Value = 5508.222
$ws.Range("I99").Value = 2012.5
$ws.Range("K99").Value = 2012.5
$ws.Range("M99").Value = -514.5
$ws.Range("H105").Value = 3573.72
$ws.Range("I105").Value = 3126.9167
$ws.Range("K105").Value = 3126.9167
$ws.Range("M105").Value = -1379.9167
$ws.Range("H107").Value = 446.14285
$ws.Range("I107").Value = 287.72726
$ws.Range("K107").Value = 287.72726
$ws.Range("M107").Value = 1632.27274
$ws.Range("H111").Value = 60000
$ws.Range("J111").Value = 60000
$ws.Range("L111").Value = 60000
$ws.Range("N111").Value = -68180
$ws.Range("H113").Value = 1650.1333
$ws.Range("I113").Value = 1286.6923
$ws.Range("K113").Value = 1286.6923
$ws.Range("M113").Value = 883.3077000000001
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").ClearContents()
$ws.Range("N114").Value = 0
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").ClearContents()
$ws.Range("N116").Value = 0
$ws.Range("H120").Value = 39998
$ws.Range("I120").Value = 39998
$ws.Range("J120").Value = 39998
$ws.Range("K120").Value = 39998
$ws.Range("L120").Value = 39998
$ws.Range("M120").Value = -47256
$ws.Range("N120").ClearContents()
$ws.Range("H126").Value = 5508.222
$ws.Range("I126").Value = 2012.5
$ws.Range("K126").Value = 6037.5
$ws.Range("M126").Value = -3567.5
$ws.Range("H136").Value = 3342.2222
$ws.Range("I136").Value = 3297.2856
$ws.Range("K136").Value = 9891.856800000001
$ws.Range("M136").Value = -7341.856800000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 22000092
$ws.Range("I4").Value = 22000092
$ws.Range("K4").Value = 66000276
$ws.Range("M4").Value = -66000164
$ws.Range("H68").Value = 15631625
$ws.Range("I68").Value = 500
$ws.Range("K68").Value = 1500
$ws.Range("M68").Value = -689
$ws.Range("H71").Value = 15631625
$ws.Range("I71").Value = 500
$ws.Range("K71").Value = 4500
$ws.Range("M71").Value = -444
$ws.Range("H94").Value = 14240.375
$ws.Range("I94").Value = 4807.6665
$ws.Range("K94").Value = 14422.9995
$ws.Range("M94").Value = -13746.9995
$ws.Range("H112").Value = 13904.833
$ws.Range("I112").Value = 1716.5
$ws.Range("K112").Value = 5149.5
$ws.Range("M112").Value = -4041.5
$ws.Range("H117").Value = 20979.4
$ws.Range("I117").Value = 1299.3334
$ws.Range("K117").Value = 3898.0002
$ws.Range("M117").Value = -456.0001999999999
$ws.Range("H121").Value = 4529.1177
$ws.Range("I121").Value = 3444.25
$ws.Range("J121").Value = 4862.923
$ws.Range("K121").Value = 10332.75
$ws.Range("L121").Value = 14588.769
$ws.Range("M121").Value = -9022.75
$ws.Range("N121").Value = -17208.769

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4357.4287
$ws.Range("I70").Value = 4357.4287
$ws.Range("K70").Value = 4357.4287
$ws.Range("M70").Value = -4087.4287
$ws.Range("H73").Value = 4357.4287
$ws.Range("I73").Value = 4357.4287
$ws.Range("K73").Value = 4357.4287
$ws.Range("M73").Value = -3421.4287
$ws.Range("H80").Value = 3108.5
$ws.Range("I80").Value = 2119.5
$ws.Range("J80").Value = 3504.1
$ws.Range("K80").Value = 2119.5
$ws.Range("L80").Value = 3504.1
$ws.Range("M80").Value = -1121.5
$ws.Range("N80").Value = -5500.1
$ws.Range("H83").Value = 3108.5
$ws.Range("I83").Value = 2119.5
$ws.Range("J83").Value = 3504.1
$ws.Range("K83").Value = 10597.5
$ws.Range("L83").Value = 17520.5
$ws.Range("M83").Value = -5605.5
$ws.Range("N83").Value = -27504.5
$ws.Range("H102").Value = 1811.0769
$ws.Range("I102").Value = 1545.3334
$ws.Range("K102").Value = 1545.3334
$ws.Range("M102").Value = 76.66660000000002
$ws.Range("H132").Value = 1633.9445
$ws.Range("I132").Value = 1700.7059
$ws.Range("J132").Value = 499
$ws.Range("K132").Value = 5102.1177
$ws.Range("L132").Value = 1497
$ws.Range("M132").Value = -2572.1177
$ws.Range("N132").Value = -6557

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4603
$ws.Range("I62").Value = 3460.25
$ws.Range("J62").Value = 5364.8335
$ws.Range("K62").Value = 3460.25
$ws.Range("L62").Value = 5364.8335
$ws.Range("M62").Value = -2836.25
$ws.Range("N62").Value = -6612.8335
$ws.Range("H65").Value = 4603
$ws.Range("I65").Value = 3460.25
$ws.Range("J65").Value = 5364.8335
$ws.Range("K65").Value = 17301.25
$ws.Range("L65").Value = 26824.1675
$ws.Range("M65").Value = -14181.25
$ws.Range("N65").Value = -33064.1675
$ws.Range("H122").Value = 867.8570999999999
$ws.Range("I122").Value = 867.8570999999999
$ws.Range("K122").Value = 2603.5713
$ws.Range("M122").Value = -153.5712999999996
